# A05c - Teacher Interview Summary: trim to a 4-line title, drop the
# "**" / BIG PICTURE lead-in, drop the rubric TOTAL row, and drop the
# "File name:" / "Submit to Canvas:" submission bullets.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Title paragraph -> 4-line title ("Teacher Interview Summary" /
#    "TCE 486/586A" / "Spring 2026 (Edwards)" / "10 points possible"),
#    each line its own run, separated by manual line breaks.
# ---------------------------------------------------------------------
$titlePara = $d.Paragraphs(1)
$titleRange = $d.Range($titlePara.Range.Start, $titlePara.Range.End)
$titleRange.Text = ""

$lines = @("Teacher Interview Summary", "TCE 486/586A", "Spring 2026 (Edwards)", "10 points possible")
$joined = [string]::Join("|", $lines)

$insertPoint = $d.Range(0, 0)
$insertPoint.InsertBefore($joined)

# Turn each "|" separator into a real manual line break (<w:br/>), scoped
# to paragraph 1 only.
$p1Range = $d.Paragraphs(1).Range
$p1Range.Find.Execute("|", $false, $false, $false, $false, $false, $true, 1, $false, "^l", 2)

# Format each line's text run (bold, black, Inter) while leaving the
# break runs untouched, which keeps the breaks as their own bare <w:r>.
$pos = 0
foreach ($line in $lines) {
    $segStart = $pos
    $segEnd = $pos + $line.Length
    $seg = $d.Range($segStart, $segEnd)
    $seg.Font.Bold = $true
    $seg.Font.Color = 0
    $seg.Font.Name = "Inter"
    $pos = $segEnd + 1
}

# ---------------------------------------------------------------------
# 2) Drop the old "**" paragraph and the "BIG PICTURE" paragraph that
#    used to follow the title.
# ---------------------------------------------------------------------
$removeStart = $d.Paragraphs(2).Range.Start
$removeEnd = $d.Paragraphs(3).Range.End
$d.Range($removeStart, $removeEnd).Delete()

# ---------------------------------------------------------------------
# 3) Remove the rubric's trailing TOTAL row.
# ---------------------------------------------------------------------
$table = $d.Tables(1)
$lastRow = $table.Rows($table.Rows.Count)
$lastRow.Delete()

# ---------------------------------------------------------------------
# 4) Submission Guidelines: keep "Format:", drop "File name:" and
#    "Submit to Canvas:" bullets.
# ---------------------------------------------------------------------
$fileNameIndex = -1
$canvasIndex = -1
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $text = $d.Paragraphs($i).Range.Text
    if ($text.StartsWith("File name:")) {
        $fileNameIndex = $i
    }
    if ($text.StartsWith("Submit to Canvas:")) {
        $canvasIndex = $i
    }
}

if ($fileNameIndex -gt 0 -and $canvasIndex -gt 0) {
    $subStart = $d.Paragraphs($fileNameIndex).Range.Start
    $subEnd = $d.Paragraphs($canvasIndex).Range.End
    $d.Range($subStart, $subEnd).Delete()
}

Write-Output "done"
